# Applies the IFRC_HIP.xlsx edit:
#  - Replace comma-separated lists with colon-separated lists in several
#    "haz_potlink" style cells (D8, D10, E12:F12, D15, D16, D18, D20).
#  - Apply the "Normal" style used throughout the sheet (style index 1)
#    to C14 and C19, which previously used the default style (index 0).
#  - Update the selection in the sheet view to match the last-edited
#    range and update the default column width very slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace commas with colons in the haz_potlink-style cells ---
$targets = @("D8", "D10", "E12", "F12", "D15", "D16", "D18", "D20")
foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $cell.Value = $cell.Value2 -replace ",", ":"
}

# --- Apply style used elsewhere in the column (matches C13/C15 etc.) ---
$ws.Range("C14").Style = $ws.Range("C13").Style
$ws.Range("C19").Style = $ws.Range("C18").Style

# --- Update the selection to match the edited cells, finishing on F1 ---
$ws.Range("D8,D10,D15:D16,D18,D20,E12:F12,F1").Select()
$ws.Range("F1").Activate()
